$wb = $excel.ActiveWorkbook

# Scheduled data refresh: update market-price derived columns (H:N)
# for the affected leve rows across each job sheet.

# ALC row 45
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value2 = 3558.4285
$ws.Range("I45").Value2 = 445
$ws.Range("J45").Value2 = 4803.8
$ws.Range("K45").Value2 = 1335
$ws.Range("L45").Value2 = 14411.4
$ws.Range("M45").Value2 = -1143
$ws.Range("N45").Value2 = -14795.4

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 9308.691999999999
$ws.Range("I62").Value2 = 1751.0834
$ws.Range("J62").Value2 = 100000
$ws.Range("K62").Value2 = 1751.0834
$ws.Range("L62").Value2 = 100000
$ws.Range("M62").Value2 = -1127.0834
$ws.Range("N62").Value2 = -101248

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 9308.691999999999
$ws.Range("I65").Value2 = 1751.0834
$ws.Range("J65").Value2 = 100000
$ws.Range("K65").Value2 = 8755.416999999999
$ws.Range("L65").Value2 = 500000
$ws.Range("M65").Value2 = -5635.416999999999
$ws.Range("N65").Value2 = -506240

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value2 = 2036.8
$ws.Range("I88").Value2 = 1921
$ws.Range("J88").Value2 = 2500
$ws.Range("K88").Value2 = 1921
$ws.Range("L88").Value2 = 2500
$ws.Range("M88").Value2 = -1515
$ws.Range("N88").Value2 = -3312

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value2 = 2036.8
$ws.Range("I91").Value2 = 1921
$ws.Range("J91").Value2 = 2500
$ws.Range("K91").Value2 = 1921
$ws.Range("L91").Value2 = 2500
$ws.Range("M91").Value2 = -517
$ws.Range("N91").Value2 = -5308

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 1876.3182
$ws.Range("I86").Value2 = 1894.238
$ws.Range("J86").Value2 = 1500
$ws.Range("K86").Value2 = 1894.238
$ws.Range("L86").Value2 = 1500
$ws.Range("M86").Value2 = -771.2380000000001
$ws.Range("N86").Value2 = -3746

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value2 = 1876.3182
$ws.Range("I89").Value2 = 1894.238
$ws.Range("J89").Value2 = 1500
$ws.Range("K89").Value2 = 9471.190000000001
$ws.Range("L89").Value2 = 7500
$ws.Range("M89").Value2 = -3855.190000000001
$ws.Range("N89").Value2 = -18732

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value2 = 3250
$ws.Range("I62").Value2 = 3250
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 3250
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = -2626

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value2 = 3250
$ws.Range("I65").Value2 = 3250
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 16250
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = -13130

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 1284.9656
$ws.Range("I107").Value2 = 1342.75
$ws.Range("J107").Value2 = 1213.8462
$ws.Range("K107").Value2 = 1342.75
$ws.Range("L107").Value2 = 1213.8462
$ws.Range("M107").Value2 = 577.25
$ws.Range("N107").Value2 = -5053.8462

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 660.1905
$ws.Range("I12").Value2 = 98.333336
$ws.Range("J12").Value2 = 1081.5834
$ws.Range("K12").Value2 = 295.000008
$ws.Range("L12").Value2 = 3244.7502
$ws.Range("M12").Value2 = -122.000008
$ws.Range("N12").Value2 = -3590.7502

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value2 = 252.6
$ws.Range("I98").Value2 = 287.66666
$ws.Range("J98").Value2 = 200
$ws.Range("K98").Value2 = 862.9999799999999
$ws.Range("L98").Value2 = 600
$ws.Range("M98").Value2 = 635.0000200000001
$ws.Range("N98").Value2 = -3596

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 686.06665
$ws.Range("I113").Value2 = 569
$ws.Range("J113").Value2 = 861.6667
$ws.Range("K113").Value2 = 1707
$ws.Range("L113").Value2 = 2585.0001
$ws.Range("M113").Value2 = 463
$ws.Range("N113").Value2 = -6925.0001

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value2 = 2960.08
$ws.Range("I117").Value2 = 409.6
$ws.Range("J117").Value2 = 3597.7
$ws.Range("K117").Value2 = 1228.8
$ws.Range("L117").Value2 = 10793.1
$ws.Range("M117").Value2 = 2213.2
$ws.Range("N117").Value2 = -17677.1

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value2 = 1452.9445
$ws.Range("I129").Value2 = 795
$ws.Range("J129").Value2 = 2768.8333
$ws.Range("K129").Value2 = 2385
$ws.Range("L129").Value2 = 8306.499899999999
$ws.Range("M129").Value2 = 2615
$ws.Range("N129").Value2 = -18306.4999

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 891.6598
$ws.Range("I131").Value2 = 408.625
$ws.Range("J131").Value2 = 935.0787
$ws.Range("K131").Value2 = 1225.875
$ws.Range("L131").Value2 = 2805.2361
$ws.Range("M131").Value2 = 3814.125
$ws.Range("N131").Value2 = -12885.2361

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2972.04
$ws.Range("I80").Value2 = 2845.25
$ws.Range("J80").Value2 = 3089.077
$ws.Range("K80").Value2 = 2845.25
$ws.Range("L80").Value2 = 3089.077
$ws.Range("M80").Value2 = -1847.25
$ws.Range("N80").Value2 = -5085.077

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value2 = 2972.04
$ws.Range("I83").Value2 = 2845.25
$ws.Range("J83").Value2 = 3089.077
$ws.Range("K83").Value2 = 14226.25
$ws.Range("L83").Value2 = 15445.385
$ws.Range("M83").Value2 = -9234.25
$ws.Range("N83").Value2 = -25429.385

# GSM row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value2 = 9856.666999999999
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 9856.666999999999
$ws.Range("K93").Value2 = 0
$ws.Range("L93").Value2 = 9856.666999999999
$ws.Range("N93").Value2 = -13600.667

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 3733.1667
$ws.Range("I122").Value2 = 5950
$ws.Range("J122").Value2 = 2624.75
$ws.Range("K122").Value2 = 17850
$ws.Range("L122").Value2 = 7874.25
$ws.Range("M122").Value2 = -15400
$ws.Range("N122").Value2 = -12774.25

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 8595.675999999999
$ws.Range("I126").Value2 = 2125
$ws.Range("J126").Value2 = 13525.714
$ws.Range("K126").Value2 = 6375
$ws.Range("L126").Value2 = 40577.142
$ws.Range("M126").Value2 = -3905
$ws.Range("N126").Value2 = -45517.142

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 5885356
$ws.Range("I7").Value2 = 9092335
$ws.Range("J7").Value2 = 5895.5
$ws.Range("K7").Value2 = 9092335
$ws.Range("L7").Value2 = 5895.5
$ws.Range("M7").Value2 = -9092223
$ws.Range("N7").Value2 = -6119.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 3465.8333
$ws.Range("I122").Value2 = 2755.5557
$ws.Range("J122").Value2 = 4176.1113
$ws.Range("K122").Value2 = 8266.667099999999
$ws.Range("L122").Value2 = 12528.3339
$ws.Range("M122").Value2 = -5816.667099999999
$ws.Range("N122").Value2 = -17428.3339

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value2 = 5885356
$ws.Range("I126").Value2 = 9092335
$ws.Range("J126").Value2 = 5895.5
$ws.Range("K126").Value2 = 27277005
$ws.Range("L126").Value2 = 17686.5
$ws.Range("M126").Value2 = -27274535
$ws.Range("N126").Value2 = -22626.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 2985.3684
$ws.Range("I132").Value2 = 1984.091
$ws.Range("J132").Value2 = 4362.125
$ws.Range("K132").Value2 = 5952.272999999999
$ws.Range("L132").Value2 = 13086.375
$ws.Range("M132").Value2 = -3422.272999999999
$ws.Range("N132").Value2 = -18146.375

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value2 = 48352.2
$ws.Range("I133").Value2 = 0
$ws.Range("J133").Value2 = 48352.2
$ws.Range("K133").Value2 = 0
$ws.Range("L133").Value2 = 48352.2
$ws.Range("N133").Value2 = -53412.2
